$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header block updates
# ---------------------------------------------------------------------------
# D2 / B7 / B11 / B13 / E13 / H15 / I15 / J15 keep the same captions, only the
# VALOR MORA total and Cant. Trabajadores count change.
$ws.Range("E11").Value = 852800
$ws.Range("C13").Value = 2

# ---------------------------------------------------------------------------
# 2. Re-order the existing debtor's (SHIRLY) period rows 16-21 so period 2411
#    (the most recent) comes first, followed by 2410, 2409, 2408, 2407, 2405.
# ---------------------------------------------------------------------------
$periodsShirly = @("2411", "2410", "2409", "2408", "2407", "2405")
$valoresShirly = @(48534, 104000, 104000, 104000, 104000, 104000)
for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodsShirly[$i]
    $ws.Range("F$row").Value = $valoresShirly[$i]
}

# ---------------------------------------------------------------------------
# 3. Insert six fresh rows (22-27) for the new debtor JAVIER ANTONIO HERRERA
#    PALMERA, copying the formatting of the SHIRLY block (rows 16-21) so the
#    borders/styles match exactly, then overwrite with the new data.
# ---------------------------------------------------------------------------
$ws.Rows("22:27").Insert()
$ws.Range("B16:J21").Copy()
$ws.Range("B22:J27").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$periodsJavier = @("2411", "2410", "2409", "2408", "2407", "2405")
$valoresJavier = @(24266, 52000, 52000, 52000, 52000, 52000)
for ($i = 0; $i -lt 6; $i++) {
    $row = 22 + $i
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1002299740"
    $ws.Range("D$row").Value = "JAVIER ANTONIO HERRERA PALMERA"
    $ws.Range("E$row").Value = $periodsJavier[$i]
    $ws.Range("F$row").Value = $valoresJavier[$i]
    $ws.Range("G$row").Value = 1300000
}

# ---------------------------------------------------------------------------
# 4. Column D must fit the longest name now in the sheet.
# ---------------------------------------------------------------------------
$ws.Columns("D").AutoFit()
